$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was row 15)
$ws.Cells.Item(2, 4).Value = 44874
$ws.Cells.Item(2, 12).Value = 'Segunda'
$ws.Cells.Item(2, 13).Value = 250
$ws.Cells.Item(2, 14).Value = 22000
$ws.Cells.Item(2, 15).Value = 23000
$ws.Cells.Item(2, 16).Value = 22500
$ws.Cells.Item(2, 19).Value = 1875

# Row 3 (was row 23)
$ws.Cells.Item(3, 4).Value = 44496
$ws.Cells.Item(3, 12).Value = 'Primera'
$ws.Cells.Item(3, 14).Value = 23000
$ws.Cells.Item(3, 15).Value = 24000
$ws.Cells.Item(3, 16).Value = 23500
$ws.Cells.Item(3, 17).Value = '$/caja 12 kilos'
$ws.Cells.Item(3, 19).Value = 1958
$ws.Cells.Item(3, 20).Value = 12

# Row 4 (was row 10)
$ws.Cells.Item(4, 4).Value = 44839
$ws.Cells.Item(4, 12).Value = 'Segunda'
$ws.Cells.Item(4, 14).Value = 26000
$ws.Cells.Item(4, 15).Value = 27000
$ws.Cells.Item(4, 16).Value = 26500
$ws.Cells.Item(4, 17).Value = '$/caja 12 kilos'
$ws.Cells.Item(4, 19).Value = 2208
$ws.Cells.Item(4, 20).Value = 12

# Row 5 (was row 9)
$ws.Cells.Item(5, 4).Value = 44489
$ws.Cells.Item(5, 12).Value = 'Primera'
$ws.Cells.Item(5, 13).Value = 200
$ws.Cells.Item(5, 14).Value = 24000
$ws.Cells.Item(5, 15).Value = 25000
$ws.Cells.Item(5, 16).Value = 24500
$ws.Cells.Item(5, 17).Value = '$/caja 12 kilos'
$ws.Cells.Item(5, 19).Value = 2042
$ws.Cells.Item(5, 20).Value = 12

# Row 6 (was row 4)
$ws.Cells.Item(6, 4).Value = 45126
$ws.Cells.Item(6, 12).Value = 'Primera'
$ws.Cells.Item(6, 13).Value = 160
$ws.Cells.Item(6, 14).Value = 14000
$ws.Cells.Item(6, 15).Value = 15000
$ws.Cells.Item(6, 16).Value = 14375
$ws.Cells.Item(6, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(6, 19).Value = 1438
$ws.Cells.Item(6, 20).Value = 10

# Row 7 (was row 5)
$ws.Cells.Item(7, 4).Value = 45126
$ws.Cells.Item(7, 12).Value = 'Segunda'
$ws.Cells.Item(7, 13).Value = 180
$ws.Cells.Item(7, 14).Value = 13000
$ws.Cells.Item(7, 15).Value = 13000
$ws.Cells.Item(7, 16).Value = 13000
$ws.Cells.Item(7, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(7, 18).Value = 'Región de Coquimbo'
$ws.Cells.Item(7, 19).Value = 1300
$ws.Cells.Item(7, 20).Value = 10

# Row 9 (was row 12)
$ws.Cells.Item(9, 4).Value = 44482
$ws.Cells.Item(9, 13).Value = 160
$ws.Cells.Item(9, 14).Value = 25000
$ws.Cells.Item(9, 15).Value = 26000
$ws.Cells.Item(9, 16).Value = 25500
$ws.Cells.Item(9, 19).Value = 2125

# Row 10 (was row 22)
$ws.Cells.Item(10, 4).Value = 44545
$ws.Cells.Item(10, 12).Value = 'Primera'
$ws.Cells.Item(10, 13).Value = 200
$ws.Cells.Item(10, 14).Value = 23000
$ws.Cells.Item(10, 15).Value = 24000
$ws.Cells.Item(10, 16).Value = 23500
$ws.Cells.Item(10, 17).Value = '$/bandeja 12 kilos'
$ws.Cells.Item(10, 19).Value = 1958

# Row 12 (was row 21)
$ws.Cells.Item(12, 4).Value = 44441
$ws.Cells.Item(12, 13).Value = 100
$ws.Cells.Item(12, 14).Value = 29000
$ws.Cells.Item(12, 15).Value = 30000
$ws.Cells.Item(12, 16).Value = 29500
$ws.Cells.Item(12, 19).Value = 2458

# Row 13 (was row 19)
$ws.Cells.Item(13, 4).Value = 45147
$ws.Cells.Item(13, 13).Value = 270
$ws.Cells.Item(13, 14).Value = 17000
$ws.Cells.Item(13, 15).Value = 18000
$ws.Cells.Item(13, 16).Value = 17500
$ws.Cells.Item(13, 17).Value = '$/caja 10 kilos'
$ws.Cells.Item(13, 19).Value = 1750
$ws.Cells.Item(13, 20).Value = 10

# Row 14 (was row 6)
$ws.Cells.Item(14, 4).Value = 44783
$ws.Cells.Item(14, 12).Value = 'Tercera'
$ws.Cells.Item(14, 14).Value = 27000
$ws.Cells.Item(14, 15).Value = 28000
$ws.Cells.Item(14, 16).Value = 27500
$ws.Cells.Item(14, 19).Value = 2292

# Row 15 (was row 3)
$ws.Cells.Item(15, 4).Value = 44160
$ws.Cells.Item(15, 13).Value = 200
$ws.Cells.Item(15, 14).Value = 19000
$ws.Cells.Item(15, 15).Value = 20000
$ws.Cells.Item(15, 16).Value = 19500
$ws.Cells.Item(15, 17).Value = '$/caja 13 kilos'
$ws.Cells.Item(15, 19).Value = 1500
$ws.Cells.Item(15, 20).Value = 13

# Row 16 (was row 7)
$ws.Cells.Item(16, 4).Value = 44860
$ws.Cells.Item(16, 12).Value = 'Primera'
$ws.Cells.Item(16, 14).Value = 23000
$ws.Cells.Item(16, 15).Value = 24000
$ws.Cells.Item(16, 16).Value = 23500
$ws.Cells.Item(16, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(16, 19).Value = 1958

# Row 17 (was row 24)
$ws.Cells.Item(17, 4).Value = 44167
$ws.Cells.Item(17, 12).Value = 'Segunda'
$ws.Cells.Item(17, 13).Value = 200
$ws.Cells.Item(17, 14).Value = 18000
$ws.Cells.Item(17, 15).Value = 19000
$ws.Cells.Item(17, 16).Value = 18500
$ws.Cells.Item(17, 17).Value = '$/caja 13 kilos'
$ws.Cells.Item(17, 19).Value = 1423
$ws.Cells.Item(17, 20).Value = 13

# Row 18 (was row 16)
$ws.Cells.Item(18, 4).Value = 44475
$ws.Cells.Item(18, 12).Value = 'Especial'
$ws.Cells.Item(18, 13).Value = 200
$ws.Cells.Item(18, 14).Value = 32000
$ws.Cells.Item(18, 15).Value = 33000
$ws.Cells.Item(18, 16).Value = 32500
$ws.Cells.Item(18, 17).Value = '$/caja 12 kilos'
$ws.Cells.Item(18, 19).Value = 2708
$ws.Cells.Item(18, 20).Value = 12

# Row 19 (was row 2)
$ws.Cells.Item(19, 4).Value = 44811
$ws.Cells.Item(19, 13).Value = 100
$ws.Cells.Item(19, 14).Value = 29000
$ws.Cells.Item(19, 15).Value = 30000
$ws.Cells.Item(19, 16).Value = 29500
$ws.Cells.Item(19, 17).Value = '$/caja 12 kilos'
$ws.Cells.Item(19, 19).Value = 2458
$ws.Cells.Item(19, 20).Value = 12

# Row 20 (was row 26)
$ws.Cells.Item(20, 4).Value = 44524
$ws.Cells.Item(20, 14).Value = 23000
$ws.Cells.Item(20, 15).Value = 24000
$ws.Cells.Item(20, 16).Value = 23500
$ws.Cells.Item(20, 17).Value = '$/caja 12 kilos'
$ws.Cells.Item(20, 19).Value = 1958
$ws.Cells.Item(20, 20).Value = 12

# Row 21 (was row 13)
$ws.Cells.Item(21, 4).Value = 44846
$ws.Cells.Item(21, 13).Value = 160
$ws.Cells.Item(21, 14).Value = 24000
$ws.Cells.Item(21, 15).Value = 25000
$ws.Cells.Item(21, 16).Value = 24500
$ws.Cells.Item(21, 19).Value = 2042

# Row 22 (was row 14)
$ws.Cells.Item(22, 4).Value = 44846
$ws.Cells.Item(22, 12).Value = 'Segunda'
$ws.Cells.Item(22, 13).Value = 100
$ws.Cells.Item(22, 14).Value = 22000
$ws.Cells.Item(22, 15).Value = 23000
$ws.Cells.Item(22, 16).Value = 22500
$ws.Cells.Item(22, 17).Value = '$/caja 12 kilos'
$ws.Cells.Item(22, 19).Value = 1875

# Row 23 (was row 20)
$ws.Cells.Item(23, 4).Value = 44468
$ws.Cells.Item(23, 14).Value = 29000
$ws.Cells.Item(23, 15).Value = 30000
$ws.Cells.Item(23, 16).Value = 29500
$ws.Cells.Item(23, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(23, 19).Value = 2950
$ws.Cells.Item(23, 20).Value = 10

# Row 24 (was row 17)
$ws.Cells.Item(24, 4).Value = 45125
$ws.Cells.Item(24, 12).Value = 'Primera'
$ws.Cells.Item(24, 13).Value = 160
$ws.Cells.Item(24, 14).Value = 14000
$ws.Cells.Item(24, 15).Value = 15000
$ws.Cells.Item(24, 16).Value = 14375
$ws.Cells.Item(24, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(24, 19).Value = 1438
$ws.Cells.Item(24, 20).Value = 10

# Row 25 (was row 18)
$ws.Cells.Item(25, 4).Value = 45125
$ws.Cells.Item(25, 13).Value = 180
$ws.Cells.Item(25, 14).Value = 13000
$ws.Cells.Item(25, 15).Value = 13000
$ws.Cells.Item(25, 16).Value = 13000
$ws.Cells.Item(25, 17).Value = '$/bandeja 10 kilos'
$ws.Cells.Item(25, 19).Value = 1300

# Row 26 (was row 25)
$ws.Cells.Item(26, 4).Value = 44776
$ws.Cells.Item(26, 12).Value = 'Segunda'
$ws.Cells.Item(26, 13).Value = 160
$ws.Cells.Item(26, 14).Value = 29000
$ws.Cells.Item(26, 15).Value = 30000
$ws.Cells.Item(26, 16).Value = 29500
$ws.Cells.Item(26, 17).Value = '$/caja 10 kilos'
$ws.Cells.Item(26, 19).Value = 2950
$ws.Cells.Item(26, 20).Value = 10
